{"js": "const pairs = [\n  [\"2024-07-30 Tuesday\", \"2024-07-31 Wednesday\"],\n  [\"67+13=\", \"5+75=\"],\n  [\"82-51=\", \"22+23=\"],\n  [\"15+1=\", \"56+12=\"],\n  [\"76-33=\", \"25+28=\"],\n  [\"62+6=\", \"59+5=\"],\n  [\"3+88=\", \"13-7=\"],\n  [\"79-40=\", \"38-1=\"],\n  [\"38+35=\", \"67-53=\"],\n  [\"64-31=\", \"64-61=\"],\n  [\"79-9=\", \"31+58=\"],\n  [\"59-28=\", \"77+3=\"],\n  [\"73-60=\", \"32+34=\"],\n  [\"89-10=\", \"82-57=\"],\n  [\"77-49=\", \"90+9=\"],\n  [\"65-27=\", \"63+8=\"],\n  [\"45+8=\", \"56+43=\"],\n  [\"48+33=\", \"49+5=\"],\n  [\"50-28=\", \"92-1=\"],\n  [\"90+4=\", \"42-32=\"],\n  [\"47+7=\", \"49+19=\"],\n  [\"91-91=\", \"98-20=\"],\n  [\"94-13=\", \"18+57=\"],\n  [\"4+45=\", \"37+17=\"],\n  [\"71-53=\", \"93-50=\"],\n  [\"64-39=\", \"6+67=\"],\n  [\"79-5=\", \"60+32=\"],\n  [\"79-39=\", \"22+40=\"],\n  [\"5+13=\", \"84+7=\"],\n  [\"76-16=\", \"39+40=\"],\n  [\"25-16=\", \"51+22=\"],\n  [\"55-34=\", \"27+44=\"],\n  [\"56+34=\", \"73-68=\"],\n  [\"7+40=\", \"97-89=\"],\n  [\"76+12=\", \"77-58=\"],\n  [\"88-28=\", \"80-61=\"],\n  [\"30-4=\", \"7+30=\"],\n  [\"70-3=\", \"15+32=\"],\n  [\"29-8=\", \"15+58=\"],\n  [\"57+22=\", \"66-48=\"],\n  [\"34+43=\", \"73-23=\"],\n  [\"1+91=\", \"66-46=\"],\n  [\"88-1=\", \"85-12=\"],\n  [\"81+14=\", \"8+69=\"],\n  [\"13+7=\", \"47+28=\"],\n  [\"50-40=\", \"23+34=\"],\n  [\"33+20=\", \"17-15=\"],\n  [\"5-3=\", \"92-69=\"],\n  [\"46+12=\", \"97-32=\"],\n  [\"88-37=\", \"37+21=\"],\n  [\"47+43=\", \"27+13=\"],\n  [\"20+5=\", \"65-26=\"],\n  [\"37+28=\", \"34+26=\"],\n  [\"23+53=\", \"15+28=\"],\n  [\"98-21=\", \"16+5=\"],\n  [\"54-47=\", \"45+50=\"],\n  [\"43-41=\", \"89-0=\"],\n  [\"4+2=\", \"40+3=\"],\n  [\"33-6=\", \"96+1=\"],\n  [\"50-34=\", \"88-71=\"],\n  [\"65+33=\", \"70-24=\"],\n  [\"68+0=\", \"52-24=\"],\n  [\"36-7=\", \"28+64=\"],\n  [\"51+11=\", \"31+17=\"],\n  [\"60+13=\", \"78-29=\"],\n  [\"78-2=\", \"7+7=\"],\n  [\"0+15=\", \"90-60=\"],\n  [\"98-93=\", \"19+43=\"],\n  [\"88-27=\", \"90-56=\"],\n  [\"13+56=\", \"2-0=\"],\n  [\"51-3=\", \"4+49=\"],\n  [\"92-20=\", \"69-19=\"],\n  [\"68-19=\", \"7+81=\"],\n  [\"30+52=\", \"71-5=\"],\n  [\"48+41=\", \"24+12=\"],\n  [\"17+36=\", \"23+17=\"],\n  [\"90-39=\", \"85+0=\"],\n  [\"22+65=\", \"37-27=\"],\n  [\"42-5=\", \"55-18=\"],\n  [\"71+28=\", \"91-14=\"],\n  [\"58+18=\", \"37+15=\"],\n  [\"53-14=\", \"53-44=\"],\n  [\"99-32=\", \"57-11=\"],\n  [\"41-14=\", \"33+60=\"],\n  [\"76-42=\", \"8-2=\"],\n  [\"53+29=\", \"13+33=\"],\n  [\"18+29=\", \"42+31=\"],\n  [\"85-69=\", \"15+80=\"],\n  [\"18+53=\", \"49+6=\"],\n  [\"83+8=\", \"42+42=\"],\n  [\"21+7=\", \"71-0=\"],\n  [\"85-42=\", \"57-56=\"],\n  [\"45-33=\", \"24+34=\"],\n  [\"48+38=\", \"34-9=\"],\n  [\"1+41=\", \"65-5=\"],\n  [\"58+26=\", \"92-24=\"],\n  [\"99-12=\", \"17+61=\"],\n  [\"11-11=\", \"32-6=\"],\n  [\"67-19=\", \"57+8=\"],\n  [\"25+53=\", \"45+40=\"],\n  [\"51+44=\", \"88-6=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2024-07-30 Tuesday', '2024-07-31 Wednesday')\n    ,@('67+13=', '5+75=')\n    ,@('82-51=', '22+23=')\n    ,@('15+1=', '56+12=')\n    ,@('76-33=', '25+28=')\n    ,@('62+6=', '59+5=')\n    ,@('3+88=', '13-7=')\n    ,@('79-40=', '38-1=')\n    ,@('38+35=', '67-53=')\n    ,@('64-31=', '64-61=')\n    ,@('79-9=', '31+58=')\n    ,@('59-28=', '77+3=')\n    ,@('73-60=', '32+34=')\n    ,@('89-10=', '82-57=')\n    ,@('77-49=', '90+9=')\n    ,@('65-27=', '63+8=')\n    ,@('45+8=', '56+43=')\n    ,@('48+33=', '49+5=')\n    ,@('50-28=', '92-1=')\n    ,@('90+4=', '42-32=')\n    ,@('47+7=', '49+19=')\n    ,@('91-91=', '98-20=')\n    ,@('94-13=', '18+57=')\n    ,@('4+45=', '37+17=')\n    ,@('71-53=', '93-50=')\n    ,@('64-39=', '6+67=')\n    ,@('79-5=', '60+32=')\n    ,@('79-39=', '22+40=')\n    ,@('5+13=', '84+7=')\n    ,@('76-16=', '39+40=')\n    ,@('25-16=', '51+22=')\n    ,@('55-34=', '27+44=')\n    ,@('56+34=', '73-68=')\n    ,@('7+40=', '97-89=')\n    ,@('76+12=', '77-58=')\n    ,@('88-28=', '80-61=')\n    ,@('30-4=', '7+30=')\n    ,@('70-3=', '15+32=')\n    ,@('29-8=', '15+58=')\n    ,@('57+22=', '66-48=')\n    ,@('34+43=', '73-23=')\n    ,@('1+91=', '66-46=')\n    ,@('88-1=', '85-12=')\n    ,@('81+14=', '8+69=')\n    ,@('13+7=', '47+28=')\n    ,@('50-40=', '23+34=')\n    ,@('33+20=', '17-15=')\n    ,@('5-3=', '92-69=')\n    ,@('46+12=', '97-32=')\n    ,@('88-37=', '37+21=')\n    ,@('47+43=', '27+13=')\n    ,@('20+5=', '65-26=')\n    ,@('37+28=', '34+26=')\n    ,@('23+53=', '15+28=')\n    ,@('98-21=', '16+5=')\n    ,@('54-47=', '45+50=')\n    ,@('43-41=', '89-0=')\n    ,@('4+2=', '40+3=')\n    ,@('33-6=', '96+1=')\n    ,@('50-34=', '88-71=')\n    ,@('65+33=', '70-24=')\n    ,@('68+0=', '52-24=')\n    ,@('36-7=', '28+64=')\n    ,@('51+11=', '31+17=')\n    ,@('60+13=', '78-29=')\n    ,@('78-2=', '7+7=')\n    ,@('0+15=', '90-60=')\n    ,@('98-93=', '19+43=')\n    ,@('88-27=', '90-56=')\n    ,@('13+56=', '2-0=')\n    ,@('51-3=', '4+49=')\n    ,@('92-20=', '69-19=')\n    ,@('68-19=', '7+81=')\n    ,@('30+52=', '71-5=')\n    ,@('48+41=', '24+12=')\n    ,@('17+36=', '23+17=')\n    ,@('90-39=', '85+0=')\n    ,@('22+65=', '37-27=')\n    ,@('42-5=', '55-18=')\n    ,@('71+28=', '91-14=')\n    ,@('58+18=', '37+15=')\n    ,@('53-14=', '53-44=')\n    ,@('99-32=', '57-11=')\n    ,@('41-14=', '33+60=')\n    ,@('76-42=', '8-2=')\n    ,@('53+29=', '13+33=')\n    ,@('18+29=', '42+31=')\n    ,@('85-69=', '15+80=')\n    ,@('18+53=', '49+6=')\n    ,@('83+8=', '42+42=')\n    ,@('21+7=', '71-0=')\n    ,@('85-42=', '57-56=')\n    ,@('45-33=', '24+34=')\n    ,@('48+38=', '34-9=')\n    ,@('1+41=', '65-5=')\n    ,@('58+26=', '92-24=')\n    ,@('99-12=', '17+61=')\n    ,@('11-11=', '32-6=')\n    ,@('67-19=', '57+8=')\n    ,@('25+53=', '45+40=')\n    ,@('51+44=', '88-6=')\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $result = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $result) {\n        Write-Output \"FAILED to replace: $old\"\n    }\n}"}
